$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values in the regression table (rows 2-4)
$ws.Range("B2").Value = "-0.958***"
$ws.Range("B3").Value = "0.211*"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "-0.011"
$ws.Range("D3").Style = "Normal"
$ws.Range("B4").Value = "12.037***"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.094"
$ws.Range("D4").Style = "Normal"

# Remove the trailing "Constant" and "r2_adj" rows (rows 5 and 6)
$ws.Range("A5:D6").Delete()
